$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (default, unstyled data cell) used to avoid leaving
# a residual NumberFormat override on cells that need to be forced to text.
$plainStyle = $ws.Range("B2").Style

$ws.Range('D2').Value = '61.944.76'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '3.425.47'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '409.21'
$ws.Range('D5').Style = $plainStyle
$ws.Range('E5').Value = '  +0.57%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '128.45'
$ws.Range('D6').Style = $plainStyle
$ws.Range('E6').Value = '  -2.78%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.632'
$ws.Range('D7').Style = $plainStyle
$ws.Range('E7').Value = '  +6.03%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.738'
$ws.Range('D9').Style = $plainStyle
$ws.Range('E9').Value = '  +6.66%  '
$ws.Range('E10').Value = '  +3.20%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '42.74'
$ws.Range('D11').Style = $plainStyle
$ws.Range('E11').Value = '  +2.00%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0000221'
$ws.Range('D12').Style = $plainStyle
$ws.Range('E12').Value = '  +48.00%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '9.14'
$ws.Range('D13').Style = $plainStyle
$ws.Range('E13').Value = '  +8.69%  '
$ws.Range('E14').Value = '  -0.16%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '21.40'
$ws.Range('D15').Style = $plainStyle
$ws.Range('E15').Value = '  +7.52%  '
$ws.Range('D16').Value = '3.963.47'
$ws.Range('E16').Value = '  -0.14%  '
$ws.Range('D17').Value = '3.435.78'
$ws.Range('E17').Value = '  -0.37%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '12.57'
$ws.Range('D18').Style = $plainStyle
$ws.Range('E18').Value = '  +8.30%  '
$ws.Range('E19').Value = '  +6.68%  '
$ws.Range('D20').Value = '61.995.92'
$ws.Range('E20').Value = '  -0.24%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '452.40'
$ws.Range('D21').Style = $plainStyle
$ws.Range('E21').Value = '  +44.94%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '92.31'
$ws.Range('D22').Style = $plainStyle
$ws.Range('E22').Value = '  +9.56%  '
$ws.Range('E23').Value = '  +1.12%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '12.98'
$ws.Range('D24').Style = $plainStyle
$ws.Range('E24').Value = '  +1.68%  '
$ws.Range('E25').Value = '  +2.57%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '32.98'
$ws.Range('D26').Style = $plainStyle
$ws.Range('E26').Value = '  +10.98%  '
$ws.Range('E27').Value = '  +7.51%  '
$ws.Range('E28').Value = '  +0.38%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.69'
$ws.Range('D29').Style = $plainStyle
$ws.Range('E29').Value = '  -2.02%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.72'
$ws.Range('D30').Style = $plainStyle
$ws.Range('E30').Value = '  -2.65%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '11.99'
$ws.Range('D31').Style = $plainStyle
$ws.Range('E31').Value = '  +5.68%  '
$ws.Range('E32').Value = '  -0.83%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '43.08'
$ws.Range('D33').Style = $plainStyle
$ws.Range('E33').Value = '  -2.22%  '
$ws.Range('E34').Value = '  -0.64%  '
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('E36').Value = '  +3.07%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '54.24'
$ws.Range('D37').Style = $plainStyle
$ws.Range('E37').Value = '  +5.15%  '
$ws.Range('E38').Value = '  +0.09%  '
$ws.Range('E39').Value = '  +1.29%  '
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.323'
$ws.Range('D40').Style = $plainStyle
$ws.Range('E40').Value = '  +2.81%  '
$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.135'
$ws.Range('D41').Style = $plainStyle
$ws.Range('E41').Value = '  +7.59%  '
$ws.Range('E42').Value = '  -3.21%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '141.92'
$ws.Range('D43').Style = $plainStyle
$ws.Range('E43').Value = '  +0.70%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.26'
$ws.Range('D44').Style = $plainStyle
$ws.Range('E44').Value = '  +8.13%  '
$ws.Range('E45').Value = '  +0.70%  '
$ws.Range('E46').Value = '  +13.28%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '16.67'
$ws.Range('D47').Style = $plainStyle
$ws.Range('E47').Value = '  -0.55%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '22.41'
$ws.Range('D48').Style = $plainStyle
$ws.Range('E48').Value = '  +4.92%  '
$ws.Range('E49').Value = '  +10.22%  '
$ws.Range('D50').Value = '3.774.28'
$ws.Range('E50').Value = '  +0.23%  '
$ws.Range('E51').Value = '  +15.45%  '
